# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sub3 = [string][char]0x2083

# Column D = Price, Column E = Volume(1h)
# Each entry: row number, new Price text (or $null to leave unchanged), new Volume text (or $null)
$updates = @(
    @{ Row = 2;  D = "62.741.41";  E = "  -0.03%  " },
    @{ Row = 3;  D = "3.049.36";   E = "  +0.16%  " },
    @{ Row = 4;  D = $null;        E = "  -0.16%  " },
    @{ Row = 5;  D = "547.50";     E = "  +2.72%  " },
    @{ Row = 6;  D = "136.25";     E = "  +1.13%  " },
    @{ Row = 7;  D = $null;        E = "  -0.08%  " },
    @{ Row = 8;  D = "3.044.80";   E = "  +0.18%  " },
    @{ Row = 9;  D = "0.496";      E = "  +1.59%  " },
    @{ Row = 10; D = $null;        E = "  -2.79%  " },
    @{ Row = 11; D = "6.14";       E = "  -0.09%  " },
    @{ Row = 12; D = "0.451";      E = "  +0.92%  " },
    @{ Row = 13; D = "35.05";      E = "  +3.31%  " },
    @{ Row = 14; D = "0.0000224";  E = "  +2.61%  " },
    @{ Row = 15; D = "3.543.60";   E = "  -0.12%  " },
    @{ Row = 16; D = "62.706.55";  E = "  -0.21%  " },
    @{ Row = 17; D = "3.055.35";   E = "  -0.03%  " },
    @{ Row = 18; D = $null;        E = "  -2.42%  " },
    @{ Row = 19; D = "6.70";       E = "  +2.44%  " },
    @{ Row = 20; D = "484.51";     E = "  +4.23%  " },
    @{ Row = 21; D = "13.37";      E = "  +1.01%  " },
    @{ Row = 22; D = "0.681";      E = "  -0.62%  " },
    @{ Row = 23; D = "7.09";       E = "  +2.72%  " },
    @{ Row = 24; D = "81.20";      E = "  +4.48%  " },
    @{ Row = 25; D = "12.19";      E = "  +1.89%  " },
    @{ Row = 26; D = "0.999";      E = "  +0.13%  " },
    @{ Row = 27; D = "2.73";       E = "  +2.70%  " },
    @{ Row = 28; D = "7.86";       E = "  +1.66%  " },
    @{ Row = 29; D = "1.00";       E = $null },
    @{ Row = 30; D = "1.95";       E = "  +5.62%  " },
    @{ Row = 31; D = "25.93";      E = "  +0.63%  " },
    @{ Row = 32; D = "1.13";       E = "  -0.99%  " },
    @{ Row = 33; D = "5.75";       E = "  +7.27%  " },
    @{ Row = 34; D = "2.38";       E = "  +4.95%  " },
    @{ Row = 35; D = "55.33";      E = "  -4.92%  " },
    @{ Row = 36; D = "5.92";       E = "  +1.13%  " },
    @{ Row = 37; D = "464.28";     E = "  +0.03%  " },
    @{ Row = 38; D = "3.193.02";   E = "  -0.58%  " },
    @{ Row = 39; D = "0.0807";     E = "  +3.17%  " },
    @{ Row = 40; D = "0.0389";     E = "  -0.48%  " },
    @{ Row = 41; D = $null;        E = "  +2.40%  " },
    @{ Row = 42; D = $null;        E = "  +2.11%  " },
    @{ Row = 43; D = $null;        E = "  -2.03%  " },
    @{ Row = 44; D = "26.59";      E = "  +7.47%  " },
    @{ Row = 45; D = $null;        E = "  -0.11%  " },
    @{ Row = 46; D = "0.246";      E = "  -0.48%  " },
    @{ Row = 47; D = $null;        E = "  +2.03%  " },
    @{ Row = 48; D = $null;        E = "  +1.06%  " },
    @{ Row = 49; D = "116.70";     E = "  -4.82%  " },
    @{ Row = 50; D = ("0.0" + $sub3 + "0495"); E = "  -2.83%  " },
    @{ Row = 51; D = $null;        E = "  +3.52%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Prefix with an apostrophe so Excel stores the price as literal text
        # (many of these look like numbers, e.g. "547.50" or "3.049.36") and
        # then reset the cell style so no stray number-format/quote-prefix
        # style sticks around on save.
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.Value = "'" + $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
